# Generate Report for Handoff
#
# The CI run that produced this workbook re-generated the localization
# status report. The only substantive data change is the "Latest Handoff
# Datetime" for the row whose Source File Name is
# "6644760e-7863-423f-8374-60342da06b65.md" on the "zh-cn" worksheet: it
# advances from 2016-09-01 16:49:20 to 2016-09-01 16:49:35.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("zh-cn")

# Column H ("Latest Handoff Datetime") is stored as plain text (not a
# numeric date), so assign the literal string to preserve the existing
# "yyyy-mm-dd HH:mm:ss"-formatted text cell exactly as Excel would when a
# user / automation re-types the same-shaped value.
$ws.Range("H5").Value = "2016-09-01 16:49:35"
